$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of columns A and B for each used row (A1:B4): the
# author/textType columns were transposed for every record.
$lastRow = 4
for ($r = 1; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $valA = $cellA.Value2
    $valB = $cellB.Value2
    $cellA.Value = $valB
    $cellB.Value = $valA
}
